# Update summary data table on Sheet1: regenerate the random-integer-derived
# columns (B=Deaths, C=Detected Infections, D=Cumulative Infections,
# E=Current Asymptomatic Infections, F=Current Infectious Carriers)
# using the new `random`-based values (replacing the old numpy-based ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    9  = @{ D = 5;  E = 4 }
    10 = @{ D = 7;  E = 6 }
    11 = @{ D = 8;  E = 7 }
    12 = @{ D = 8;  E = 7 }
    14 = @{ D = 10; E = 5;  F = 4 }
    15 = @{ D = 14; E = 7;  F = 6 }
    16 = @{ D = 20; E = 12; F = 7 }
    17 = @{ D = 29; E = 21; F = 7 }
    18 = @{ D = 34; E = 26 }
    19 = @{ C = 5;  D = 38; E = 28; F = 5 }
    20 = @{ C = 7;  D = 42; E = 28; F = 7 }
    21 = @{ C = 8;  D = 51; E = 35 }
    22 = @{ C = 8;  D = 59; E = 33; F = 18 }
    23 = @{ B = 1;  D = 73; E = 41; F = 24 }
    24 = @{ B = 1;  C = 10; D = 83;  E = 49;  F = 24 }
    25 = @{ B = 1;  C = 14; D = 98;  E = 58;  F = 26 }
    26 = @{ B = 1;  C = 16; D = 118; E = 73;  F = 29 }
    27 = @{ B = 1;  C = 26; D = 141; E = 86;  F = 29 }
    28 = @{ B = 1;  C = 32; D = 165; E = 100; F = 33 }
    29 = @{ B = 2;  C = 34; D = 191; E = 115; F = 42 }
    30 = @{ B = 2;  C = 40; D = 228; E = 137; F = 51 }
    31 = @{ B = 2;  C = 45; D = 272; E = 165; F = 62 }
    32 = @{ B = 2;  C = 55; D = 326; E = 200; F = 71 }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $rowData[$col]
    }
}
